# Refresh the "Created Room ID" values (column L) in the booking sheet with
# the latest IDs captured from the most recent test run (see commit message:
# "docs - latest version used in the yaml file"). Row 64 keeps its previous
# value because it was not part of the refreshed batch.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "L2" = "220930916"
    "L9" = "220973741"
    "L10" = "220975277"
    "L11" = "220976305"
    "L12" = "220977462"
    "L13" = "220978293"
    "L14" = "220979139"
    "L15" = "220979975"
    "L16" = "220981368"
    "L17" = "220982239"
    "L18" = "220983646"
    "L19" = "220984585"
    "L20" = "220985453"
    "L21" = "220986403"
    "L22" = "220987351"
    "L23" = "220988325"
    "L24" = "220989709"
    "L25" = "220990610"
    "L26" = "220991836"
    "L27" = "220992806"
    "L28" = "220994313"
    "L29" = "220995468"
    "L30" = "220996679"
    "L31" = "220998173"
    "L32" = "220999186"
    "L33" = "221000230"
    "L34" = "221001231"
    "L35" = "221002261"
    "L36" = "221003626"
    "L37" = "221004891"
    "L38" = "221006205"
    "L39" = "221007307"
    "L40" = "221008377"
    "L41" = "221009595"
    "L42" = "221010878"
    "L43" = "221011972"
    "L44" = "221013481"
    "L45" = "221014725"
    "L46" = "221015875"
    "L47" = "221017508"
    "L48" = "221018693"
    "L49" = "221019896"
    "L50" = "221021129"
    "L51" = "221022650"
    "L52" = "221023854"
    "L53" = "221025211"
    "L54" = "221026505"
    "L55" = "221028022"
    "L56" = "221029891"
    "L57" = "221032144"
    "L58" = "221034744"
    "L59" = "221036395"
    "L60" = "221039201"
    "L61" = "221042579"
    "L62" = "221045176"
    "L63" = "221050316"
    "L65" = "221055188"
    "L66" = "221058023"
    "L67" = "221059818"
    "L68" = "221061587"
    "L69" = "221063270"
    "L70" = "221065001"
    "L71" = "221067352"
    "L72" = "221069911"
    "L73" = "221073130"
    "L74" = "221076549"
    "L75" = "221079281"
    "L76" = "221082353"
    "L77" = "221085186"
}

# The "Created Room ID" values are numeric-looking but must stay stored as
# text (shared strings), matching their original type. Writing the digits
# directly via .Value gets auto-detected as a number, so instead build each
# value as a text-formula result in a scratch cell and paste only the value
# back - this keeps the text type without touching number formats/styles.
$scratch = $ws.Range("ZZ1")
foreach ($cellRef in $newValues.Keys) {
    $scratch.Formula = '="' + $newValues[$cellRef] + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}
$scratch.ClearContents()

